$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.224.52"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "1.604.75"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3772"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.49"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3626"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08146"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.604"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.355"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001247"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "1.603.85"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06938"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.537"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("D24").Value = "23.220.93"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.438"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.058"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.294"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.409"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.777"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("D33").Value = "1.780.04"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9577"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02771"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.07445"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2519"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.126"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08782"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.407"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7103"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6537"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.333"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.011"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "133.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.200"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.14%  "
